$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 26
$ws.Range("B3").Value = 910000
$ws.Range("B4").Value = 13928571.42857143
$ws.Range("B15").Value = 9285714.285714285
$ws.Range("B26").Value = 13928571.42857143
$ws.Range("B35").Value = 10426571.42857143
$ws.Range("B36").Value = 9285714.285714285
$ws.Range("B37").Value = 13928571.42857143
$ws.Range("B38").Value = 33640857.14285715

$wb.Save()
